$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-4 and overwrite row 5, then add new row 6,
# matching the new schedule data (columns A:J).
# Columns: A=trialTrain, B=x_fixStart, C=y_fixStart, D=x_corrSteps, E=y_corrSteps,
#          F=x_nrSteps, G=y_nrSteps, H=alienID, I=praclen, J=version

$data = @(
    @(1, 2, 1, 2, 4, 0, 3, 64, 5, "train_dim1_1"),
    @(2, 1, 0, 1, 0, 0, 0, 61, 5, "train_dim1_1"),
    @(3, 4, 3, 4, 7, 0, 4, 65, 5, "train_dim1_1"),
    @(4, 1, 0, 1, 5, 0, 5, 66, 5, "train_dim1_1"),
    @(5, 3, 1, 3, 2, 0, 1, 62, 5, "train_dim1_1")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $rowData[$c]
    }
}

$ws.Range("G11").Select()
